$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 9).Value = "aa"
$ws.Cells.Item(2, 10).Value = "Agree/Accept"
$ws.Cells.Item(6, 9).Value = "sv"
$ws.Cells.Item(6, 10).Value = "Statement-opinion"
$ws.Cells.Item(15, 9).Value = "sv"
$ws.Cells.Item(15, 10).Value = "Statement-opinion"
$ws.Cells.Item(18, 9).Value = "ba"
$ws.Cells.Item(18, 10).Value = "Appreciation"
$ws.Cells.Item(22, 9).Value = "aa"
$ws.Cells.Item(22, 10).Value = "Agree/Accept"
$ws.Cells.Item(32, 9).Value = "sv"
$ws.Cells.Item(32, 10).Value = "Statement-opinion"
$ws.Cells.Item(36, 9).Value = "sv"
$ws.Cells.Item(36, 10).Value = "Statement-opinion"
$ws.Cells.Item(43, 9).Value = "ba"
$ws.Cells.Item(43, 10).Value = "Appreciation"
$ws.Cells.Item(53, 9).Value = "sv"
$ws.Cells.Item(53, 10).Value = "Statement-opinion"
$ws.Cells.Item(59, 9).Value = "qy"
$ws.Cells.Item(59, 10).Value = "Yes-No-Question"
$ws.Cells.Item(71, 9).Value = "sv"
$ws.Cells.Item(71, 10).Value = "Statement-opinion"
$ws.Cells.Item(76, 9).Value = "sd"
$ws.Cells.Item(76, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(82, 9).Value = "sv"
$ws.Cells.Item(82, 10).Value = "Statement-opinion"
$ws.Cells.Item(83, 9).Value = "aa"
$ws.Cells.Item(83, 10).Value = "Agree/Accept"
$ws.Cells.Item(85, 9).Value = "qy"
$ws.Cells.Item(85, 10).Value = "Yes-No-Question"
$ws.Cells.Item(98, 9).Value = "sv"
$ws.Cells.Item(98, 10).Value = "Statement-opinion"
$ws.Cells.Item(102, 9).Value = "sd"
$ws.Cells.Item(102, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(103, 9).Value = "sv"
$ws.Cells.Item(103, 10).Value = "Statement-opinion"
$ws.Cells.Item(118, 9).Value = "sv"
$ws.Cells.Item(118, 10).Value = "Statement-opinion"
$ws.Cells.Item(119, 9).Value = "sv"
$ws.Cells.Item(119, 10).Value = "Statement-opinion"
$ws.Cells.Item(121, 9).Value = "aa"
$ws.Cells.Item(121, 10).Value = "Agree/Accept"
$ws.Cells.Item(131, 9).Value = "sd"
$ws.Cells.Item(131, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(138, 9).Value = "sd"
$ws.Cells.Item(138, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(157, 9).Value = "aa"
$ws.Cells.Item(157, 10).Value = "Agree/Accept"
$ws.Cells.Item(158, 9).Value = "aa"
$ws.Cells.Item(158, 10).Value = "Agree/Accept"
$ws.Cells.Item(167, 9).Value = "ba"
$ws.Cells.Item(167, 10).Value = "Appreciation"
$ws.Cells.Item(191, 9).Value = "sd"
$ws.Cells.Item(191, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(238, 9).Value = "sd"
$ws.Cells.Item(238, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(242, 9).Value = "sd"
$ws.Cells.Item(242, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(248, 9).Value = "qy"
$ws.Cells.Item(248, 10).Value = "Yes-No-Question"
$ws.Cells.Item(253, 9).Value = "ba"
$ws.Cells.Item(253, 10).Value = "Appreciation"
$ws.Cells.Item(269, 9).Value = "sd"
$ws.Cells.Item(269, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(274, 9).Value = "sd"
$ws.Cells.Item(274, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(291, 9).Value = "sd"
$ws.Cells.Item(291, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(296, 9).Value = "ba"
$ws.Cells.Item(296, 10).Value = "Appreciation"
$ws.Cells.Item(298, 9).Value = "%"
$ws.Cells.Item(298, 10).Value = "Uninterpretable"
$ws.Cells.Item(300, 9).Value = "%"
$ws.Cells.Item(300, 10).Value = "Uninterpretable"
$ws.Cells.Item(301, 9).Value = "sv"
$ws.Cells.Item(301, 10).Value = "Statement-opinion"
$ws.Cells.Item(302, 9).Value = "ba"
$ws.Cells.Item(302, 10).Value = "Appreciation"
$ws.Cells.Item(309, 9).Value = "ba"
$ws.Cells.Item(309, 10).Value = "Appreciation"
$ws.Cells.Item(314, 9).Value = "sd"
$ws.Cells.Item(314, 10).Value = "Statement-non-opinion"
$ws.Cells.Item(315, 9).Value = "sv"
$ws.Cells.Item(315, 10).Value = "Statement-opinion"
$ws.Cells.Item(317, 9).Value = "sv"
$ws.Cells.Item(317, 10).Value = "Statement-opinion"
